# Auto-generated Excel COM-interop script
# Applies updated market-board snapshot values (currentAveragePrice* / LevePrice* / LeveProfit*)
# to the per-job profit sheets, matching the scheduled runner's data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2357.1794
$ws.Range("I15").Value = 2357.1794
$ws.Range("K15").Value = 7071.5382
$ws.Range("M15").Value = -6902.5382
$ws.Range("H29").Value = 2014.5
$ws.Range("J29").Value = 2909.375
$ws.Range("L29").Value = 8728.125
$ws.Range("N29").Value = -9290.125
$ws.Range("H43").Value = 4318.9
$ws.Range("I43").Value = 2301
$ws.Range("J43").Value = 5183.7144
$ws.Range("K43").Value = 2301
$ws.Range("L43").Value = 5183.7144
$ws.Range("M43").Value = -2232
$ws.Range("N43").Value = -5321.7144
$ws.Range("H52").Value = 566.6667
$ws.Range("I52").Value = 566.6667
$ws.Range("K52").Value = 1700.0001
$ws.Range("M52").Value = -1540.0001
$ws.Range("H53").Value = 1545
$ws.Range("I53").Value = 324.5
$ws.Range("J53").Value = 2358.6667
$ws.Range("K53").Value = 324.5
$ws.Range("L53").Value = 2358.6667
$ws.Range("M53").Value = 312.5
$ws.Range("N53").Value = -3632.6667
$ws.Range("H86").Value = 7898.75
$ws.Range("J86").Value = 7933.3335
$ws.Range("L86").Value = 7933.3335
$ws.Range("N86").Value = -10179.3335
$ws.Range("H87").Value = 80000
$ws.Range("J87").Value = 80000
$ws.Range("L87").Value = 80000
$ws.Range("N87").Value = -82496
$ws.Range("H89").Value = 7898.75
$ws.Range("J89").Value = 7933.3335
$ws.Range("L89").Value = 39666.6675
$ws.Range("N89").Value = -50898.6675
$ws.Range("H90").Value = 80000
$ws.Range("J90").Value = 80000
$ws.Range("L90").Value = 240000
$ws.Range("N90").Value = -252480
$ws.Range("H99").Value = 1134.2
$ws.Range("I99").Value = 1134.2
$ws.Range("K99").Value = 3402.6
$ws.Range("M99").Value = -1904.6
$ws.Range("H101").Value = 1147.7
$ws.Range("I101").Value = 685.375
$ws.Range("J101").Value = 2997
$ws.Range("K101").Value = 2056.125
$ws.Range("L101").Value = 8991
$ws.Range("M101").Value = -434.125
$ws.Range("N101").Value = -12235
$ws.Range("H112").Value = 1886.7778
$ws.Range("I112").Value = 864.8
$ws.Range("J112").Value = 2279.8462
$ws.Range("K112").Value = 2594.4
$ws.Range("L112").Value = 6839.5386
$ws.Range("M112").Value = -1486.4
$ws.Range("N112").Value = -9055.5386
$ws.Range("H116").Value = 11667.5
$ws.Range("I116").Value = 13499.5
$ws.Range("K116").Value = 13499.5
$ws.Range("M116").Value = -10057.5
$ws.Range("H118").Value = 900.4666999999999
$ws.Range("I118").Value = 231.71428
$ws.Range("K118").Value = 695.14284
$ws.Range("M118").Value = 961.85716
$ws.Range("H132").Value = 1782.4073
$ws.Range("I132").Value = 1539.6957
$ws.Range("J132").Value = 3178
$ws.Range("K132").Value = 4619.0871
$ws.Range("L132").Value = 9534
$ws.Range("M132").Value = -2089.0871
$ws.Range("N132").Value = -14594
$ws.Range("H137").Value = 2443.1667
$ws.Range("I137").Value = 1725.7142
$ws.Range("J137").Value = 3160.6191
$ws.Range("K137").Value = 5177.142599999999
$ws.Range("L137").Value = 9481.8573
$ws.Range("M137").Value = -2627.142599999999
$ws.Range("N137").Value = -14581.8573
$ws.Range("H138").Value = 3299.9375
$ws.Range("J138").Value = 4046.926
$ws.Range("L138").Value = 12140.778
$ws.Range("N138").Value = -22420.778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 279000
$ws.Range("I34").Value = 168500
$ws.Range("K34").Value = 168500
$ws.Range("M34").Value = -168229
$ws.Range("H61").Value = 8206.473
$ws.Range("I61").Value = 6262.2856
$ws.Range("J61").Value = 15011.125
$ws.Range("K61").Value = 6262.2856
$ws.Range("L61").Value = 15011.125
$ws.Range("M61").Value = -6050.2856
$ws.Range("N61").Value = -15435.125
$ws.Range("H74").Value = 20835512
$ws.Range("I74").Value = 27780000
$ws.Range("J74").Value = 2049.75
$ws.Range("K74").Value = 27780000
$ws.Range("L74").Value = 2049.75
$ws.Range("M74").Value = -27779126
$ws.Range("N74").Value = -3797.75
$ws.Range("H77").Value = 20835512
$ws.Range("I77").Value = 27780000
$ws.Range("J77").Value = 2049.75
$ws.Range("K77").Value = 138900000
$ws.Range("L77").Value = 10248.75
$ws.Range("M77").Value = -138895632
$ws.Range("N77").Value = -18984.75
$ws.Range("H136").Value = 8206.473
$ws.Range("I136").Value = 6262.2856
$ws.Range("J136").Value = 15011.125
$ws.Range("K136").Value = 18786.8568
$ws.Range("L136").Value = 45033.375
$ws.Range("M136").Value = -16236.8568
$ws.Range("N136").Value = -50133.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30362.41
$ws.Range("I31").Value = 2896.1724
$ws.Range("J31").Value = 110014.5
$ws.Range("K31").Value = 2896.1724
$ws.Range("L31").Value = 110014.5
$ws.Range("M31").Value = -2601.1724
$ws.Range("N31").Value = -110604.5
$ws.Range("H34").Value = 30362.41
$ws.Range("I34").Value = 2896.1724
$ws.Range("J34").Value = 110014.5
$ws.Range("K34").Value = 2896.1724
$ws.Range("L34").Value = 110014.5
$ws.Range("M34").Value = -2694.1724
$ws.Range("N34").Value = -110418.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 76.40000000000001
$ws.Range("J114").Value = 33
$ws.Range("L114").Value = 99
$ws.Range("N114").Value = -6607

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 14171.667
$ws.Range("I35").Value = 23015
$ws.Range("K35").Value = 23015
$ws.Range("M35").Value = -22717
$ws.Range("H122").Value = 6061.1816
$ws.Range("I122").Value = 5520.625
$ws.Range("J122").Value = 7502.6665
$ws.Range("K122").Value = 16561.875
$ws.Range("L122").Value = 22507.9995
$ws.Range("M122").Value = -14111.875
$ws.Range("N122").Value = -27407.9995
$ws.Range("H132").Value = 6869.933
$ws.Range("I132").Value = 3732.1
$ws.Range("J132").Value = 13145.6
$ws.Range("K132").Value = 11196.3
$ws.Range("L132").Value = 39436.8
$ws.Range("M132").Value = -8666.299999999999
$ws.Range("N132").Value = -44496.8
$ws.Range("H135").Value = 69996
$ws.Range("J135").Value = 69996
$ws.Range("L135").Value = 69996
$ws.Range("N135").Value = -80136

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2857.3572
$ws.Range("I16").Value = 2857.3572
$ws.Range("K16").Value = 2857.3572
$ws.Range("M16").Value = -2687.3572
$ws.Range("H132").Value = 3138.0527
$ws.Range("I132").Value = 1734.4242
$ws.Range("K132").Value = 5203.2726
$ws.Range("M132").Value = -2673.2726

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 49993.668
$ws.Range("J93").Value = 49993.668
$ws.Range("L93").Value = 49993.668
$ws.Range("N93").Value = -54985.668
$ws.Range("H100").Value = 1331.6666
$ws.Range("I100").Value = 997.5
$ws.Range("K100").Value = 1995
$ws.Range("M100").Value = -1454
